$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Coinranking.com price/volume snapshot refresh (GitHub Actions cron).
# Price/Volume(1h) cells are free-text (e.g. "1.908.52", "  +2.03%  ") rather than
# numbers, so each updated cell is forced to Text format before the new string is
# written -- otherwise Excel would auto-coerce look-alikes such as "41.90" or
# "123.00" into the numbers 41.9 / 123, dropping the trailing zero.

# --- Price (D) and Volume(1h) (E) updates for rows 2-41, 44-51 ---
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "28.264.22"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +0.16%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.912.33"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +2.22%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9995"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  -0.13%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "314.01"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +0.80%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.9989"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -0.15%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5079"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +0.17%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3933"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  -0.10%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.09357"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  -3.73%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.141"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "41.90"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +2.44%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "6.401"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  -1.86%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "20.93"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  -0.28%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "1.909.20"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +2.05%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.311"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  -1.76%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.9994"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  -0.14%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.00001125"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  -0.43%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "92.79"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  -0.28%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06586"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "17.95"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  +2.19%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.9999"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  +0.00%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.210"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  +0.75%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "28.314.41"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  +0.13%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "11.42"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.305"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  +0.93%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.594"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +1.35%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.136.74"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  +2.49%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "21.06"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  -0.85%  "
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  -0.34%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "127.44"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -0.23%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.102"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  +3.08%  "
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +1.04%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.643"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -0.04%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.611"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  -0.32%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "9.688"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  +1.24%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.06675"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  -0.84%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.02417"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.255"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  +1.22%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.2187"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  -0.23%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.267"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  +7.16%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.6416"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +0.52%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.9985"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -0.14%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "13.31"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -1.55%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.6016"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  -0.24%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.716"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  +1.46%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.275"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +1.30%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.023"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  +1.17%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "123.00"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  -0.96%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.188"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  -0.80%  "

# --- Rows 42 and 43 swap (Aptos <-> InternetComputer(DFINITY)) with new price/volume data ---
$ws.Range("B42").Value = "InternetComputer(DFINITY)"
$ws.Range("C42").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "5.011"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  +0.55%  "

$ws.Range("B43").Value = "Aptos"
$ws.Range("C43").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "11.50"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -0.24%  "
